# Apply the edits described by the commit "plano de atendimento impressão - OK"
# to the OAC (Ordem de Atendimento ao Cliente) print template.

$d = $word.ActiveDocument

# 1. The leading empty paragraph only held a leftover "_GoBack" bookmark
#    (an artifact Word drops once you've made further edits/re-saved the
#    file). Delete it so the paragraph becomes truly empty.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Company address line: append the house-number placeholder after
#    "street" -> "Endereço: street, ${numero}".
$rng = $d.Content
$foundStreet = $rng.Find.Execute("street", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundStreet) {
    $rng.Collapse(0)
    $rng.InsertAfter(", `${numero}")
}

# 3. Bump the request-received date one day forward.
$d.Content.Find.Execute("16-12-2019", $false, $false, $false, $false, $false, $true, 1, $false, "17-12-2019", 2) | Out-Null
